$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new value would otherwise be
# auto-converted to a number by Excel (losing exact text representation).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Update Price (D) and Volume/1h change (E) columns with latest scraped data
$ws.Range("D2").Value = "39.615.80"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "2.215.25"
$ws.Range("E3").Value = "  -5.31%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "298.60"
$ws.Range("E5").Value = "  -3.94%  "
$ws.Range("D6").Value = "83.69"
$ws.Range("E6").Value = "  -2.11%  "
$ws.Range("E7").Value = "  -3.02%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "0.466"
$ws.Range("E9").Value = "  -3.84%  "
$ws.Range("D10").Value = "0.0780"
$ws.Range("E10").Value = "  -3.94%  "
$ws.Range("D11").Value = "29.60"
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("D12").Value = "46.18"
$ws.Range("E12").Value = "  -11.94%  "
$ws.Range("E13").Value = "  -2.12%  "
$ws.Range("D14").Value = "2.557.78"
$ws.Range("E14").Value = "  -5.16%  "
$ws.Range("D15").Value = "6.29"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("E16").Value = "  -4.61%  "
$ws.Range("D17").Value = "2.219.89"
$ws.Range("E17").Value = "  -4.99%  "
$ws.Range("D18").Value = "0.718"
$ws.Range("E18").Value = "  -5.28%  "
$ws.Range("D19").Value = "39.543.68"
$ws.Range("E19").Value = "  -1.19%  "
$ws.Range("E20").Value = "  -2.89%  "
$ws.Range("E21").Value = "  -6.29%  "
$ws.Range("D22").Value = "64.98"
$ws.Range("E22").Value = "  -4.37%  "
$ws.Range("D23").Value = "10.40"
$ws.Range("E23").Value = "  -2.52%  "
$ws.Range("D24").Value = "232.30"
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("E26").Value = "  -5.00%  "
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("D28").Value = "22.71"
$ws.Range("E28").Value = "  -2.73%  "
$ws.Range("E29").Value = "  +2.42%  "
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("D31").Value = "32.30"
$ws.Range("E31").Value = "  -7.37%  "
$ws.Range("D32").Value = "149.30"
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("E34").Value = "  -5.60%  "
$ws.Range("E35").Value = "  -2.35%  "
$ws.Range("E36").Value = "  -2.54%  "
$ws.Range("D37").Value = "16.11"
$ws.Range("E37").Value = "  +3.15%  "
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").Value = "2.66"
$ws.Range("E40").Value = "  -5.36%  "
$ws.Range("E41").Value = "  -4.13%  "
$ws.Range("D42").Value = "3.65"
$ws.Range("E42").Value = "  -6.20%  "
$ws.Range("D43").Value = "1.926.61"
$ws.Range("E44").Value = "  -3.40%  "
$ws.Range("D45").Value = "0.0266"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("E48").Value = "  -3.93%  "
$ws.Range("D49").Value = "2.430.24"
$ws.Range("E49").Value = "  -4.69%  "
$ws.Range("D50").Value = "70.75"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").Value = "88.66"
$ws.Range("E51").Value = "  -4.49%  "

# Rows 46 and 47 changed rank order (re-sorted by the scraper); update
# coin name/link/price/volume so row contents match the new ranking.
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "9.22"
$ws.Range("E46").Value = "  -2.10%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "16.41"
$ws.Range("E47").Value = "  -6.89%  "

